# "Switching data to another sheet to test"
#
# The job-tracking table currently lives on the "jobs" tab. Move it onto
# the "notes" tab instead, leaving "jobs" empty.

$wb    = $excel.ActiveWorkbook
$jobs  = $wb.Worksheets.Item("jobs")
$notes = $wb.Worksheets.Item("notes")

$rows = 35
$cols = 9
$srcRange = $jobs.Range("A1:I35")

# Copy the whole table (values, formats, number formats, etc.) over to the
# "notes" sheet.
$srcRange.Copy($notes.Range("A1"))

# Excel's Copy materializes an empty <c> placeholder for every blank cell
# inside the rectangular copy range. Blank out those same cells on the
# destination so "notes" ends up with exactly the same (sparse) cell layout
# the table had back on "jobs".
for ($r = 1; $r -le $rows; $r++) {
  for ($c = 1; $c -le $cols; $c++) {
    if ($jobs.Cells.Item($r, $c).Value2 -eq $null) {
      $notes.Cells.Item($r, $c).Value = ""
    }
  }
}

# Now clear the source table off of "jobs" - it moved to "notes". Also
# drop the header row's bold/"customFormat" row styling so the sheet goes
# back to a plain, unformatted blank sheet.
$srcRange.Clear()
$jobs.Rows.Item(1).ClearFormats()

# "notes" is now the working sheet: select the full-sheet (row) style
# selection on it, the way you'd leave it after pasting and selecting
# everything to check the paste.
$notes.Activate()
$notes.Range("A1:XFD1048576").Select()

# Flip back to "jobs" (still the selected tab) with the view scrolled down
# and the cursor left on D37, where editing continued after the cut.
$jobs.Activate()
$jobs.Application.ActiveWindow.ScrollRow = 16
$jobs.Range("D37").Select()
